$d = $word.ActiveDocument

$replacements = @(
    @{old="652×6="; new="744×6="},
    @{old="454×7="; new="712×7="},
    @{old="435×3="; new="564×8="},
    @{old="487×7="; new="824×2="},
    @{old="990×6="; new="975×8="},
    @{old="442×4="; new="789×3="},
    @{old="579×7="; new="624×6="},
    @{old="819×3="; new="535×9="},
    @{old="683×2="; new="638×6="},
    @{old="892×5="; new="449×9="},
    @{old="372×7="; new="316×9="},
    @{old="988×6="; new="190×6="},
    @{old="581×6="; new="416×2="},
    @{old="153×5="; new="925×8="},
    @{old="837×4="; new="617×4="},
    @{old="799×4="; new="482×3="},
    @{old="502×3="; new="625×6="},
    @{old="474×4="; new="462×6="},
    @{old="725×3="; new="978×4="},
    @{old="198×6="; new="797×6="},
    @{old="980×8="; new="169×9="},
    @{old="177×6="; new="500×4="},
    @{old="994×8="; new="436×6="},
    @{old="999×2="; new="206×4="},
    @{old="605×7="; new="530×6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
